$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) updates to column F ("想去人数" / want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(4, 6).Value = 3710
$ws1.Cells.Item(5, 6).Value = 3710
$ws1.Cells.Item(6, 6).Value = 283
$ws1.Cells.Item(7, 6).Value = 5242
$ws1.Cells.Item(8, 6).Value = 570
$ws1.Cells.Item(9, 6).Value = 401
$ws1.Cells.Item(10, 6).Value = 217
$ws1.Cells.Item(11, 6).Value = 1028
$ws1.Cells.Item(13, 6).Value = 124
$ws1.Cells.Item(14, 6).Value = 44
$ws1.Cells.Item(15, 6).Value = 721
$ws1.Cells.Item(16, 6).Value = 347
$ws1.Cells.Item(19, 6).Value = 167
$ws1.Cells.Item(22, 6).Value = 6002
$ws1.Cells.Item(26, 6).Value = 6308
$ws1.Cells.Item(28, 6).Value = 21
$ws1.Cells.Item(29, 6).Value = 3245
$ws1.Cells.Item(30, 6).Value = 361
$ws1.Cells.Item(31, 6).Value = 737
$ws1.Cells.Item(32, 6).Value = 4452
$ws1.Cells.Item(34, 6).Value = 132
$ws1.Cells.Item(36, 6).Value = 1104
$ws1.Cells.Item(37, 6).Value = 95
$ws1.Cells.Item(39, 6).Value = 4
$ws1.Cells.Item(40, 6).Value = 907
$ws1.Cells.Item(41, 6).Value = 1090
$ws1.Cells.Item(42, 6).Value = 2048
$ws1.Cells.Item(43, 6).Value = 3

# Sheet "全部类型" (all types) updates to column F ("想去人数" / want-to-go count)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(7, 6).Value = 3710
$ws4.Cells.Item(8, 6).Value = 3710
$ws4.Cells.Item(9, 6).Value = 283
$ws4.Cells.Item(10, 6).Value = 5242
$ws4.Cells.Item(11, 6).Value = 570
$ws4.Cells.Item(12, 6).Value = 401
$ws4.Cells.Item(13, 6).Value = 217
$ws4.Cells.Item(14, 6).Value = 1028
$ws4.Cells.Item(16, 6).Value = 124
$ws4.Cells.Item(17, 6).Value = 44
$ws4.Cells.Item(18, 6).Value = 721
$ws4.Cells.Item(19, 6).Value = 347
$ws4.Cells.Item(23, 6).Value = 167
$ws4.Cells.Item(26, 6).Value = 6002
$ws4.Cells.Item(30, 6).Value = 6308
$ws4.Cells.Item(32, 6).Value = 21
$ws4.Cells.Item(33, 6).Value = 3245
$ws4.Cells.Item(34, 6).Value = 361
$ws4.Cells.Item(35, 6).Value = 737
$ws4.Cells.Item(36, 6).Value = 4452
$ws4.Cells.Item(39, 6).Value = 132
$ws4.Cells.Item(41, 6).Value = 1104
$ws4.Cells.Item(42, 6).Value = 95
$ws4.Cells.Item(44, 6).Value = 4
$ws4.Cells.Item(45, 6).Value = 907
$ws4.Cells.Item(46, 6).Value = 1090
$ws4.Cells.Item(48, 6).Value = 2048
$ws4.Cells.Item(49, 6).Value = 3
